$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()

# Give the newly inserted column roughly the same width as its neighbour
# (column M, "In Advance"); the engine snaps column widths to a 1/6-character
# grid, so this lands on the closest achievable width.
$ws.Columns("N:N").ColumnWidth = 9.8

# --- "Transactions" sheet used to be the active tab; make "Repayment
#     schedule" the active tab/selection instead ---
$ws.Range("R13").Select()
